$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.819.17"
$ws.Range("E2").Value = "  +6.80%  "

$ws.Range("D3").Value = "2.622.49"
$ws.Range("E3").Value = "  +9.34%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'509.80"
$ws.Range("E5").Value = "  +4.87%  "

$ws.Range("D6").Value = "'157.69"
$ws.Range("E6").Value = "  +2.52%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.604"
$ws.Range("E8").Value = "  -0.33%  "

$ws.Range("D9").Value = "2.660.94"
$ws.Range("E9").Value = "  +10.28%  "

$ws.Range("D10").Value = "'6.40"
$ws.Range("E10").Value = "  +2.00%  "

$ws.Range("E11").Value = "  +5.69%  "

$ws.Range("E12").Value = "  +3.75%  "

$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("D14").Value = "3.094.71"
$ws.Range("E14").Value = "  +9.68%  "

$ws.Range("D15").Value = "60.793.30"
$ws.Range("E15").Value = "  +6.79%  "

$ws.Range("D16").Value = "'21.83"
$ws.Range("E16").Value = "  +5.99%  "

$ws.Range("E17").Value = "  +6.07%  "

$ws.Range("D18").Value = "2.652.59"
$ws.Range("E18").Value = "  +9.99%  "

$ws.Range("D19").Value = "'4.81"
$ws.Range("E19").Value = "  +1.86%  "

$ws.Range("D20").Value = "'348.89"
$ws.Range("E20").Value = "  +7.60%  "

$ws.Range("D21").Value = "'10.53"
$ws.Range("E21").Value = "  +6.41%  "

$ws.Range("E22").Value = "  +4.75%  "

$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").Value = "'60.58"
$ws.Range("E24").Value = "  +4.40%  "

$ws.Range("E25").Value = "  +5.21%  "

$ws.Range("D26").Value = "2.759.49"
$ws.Range("E26").Value = "  +9.87%  "

$ws.Range("E27").Value = "  +4.43%  "

$ws.Range("D28").Value = "'0.985"
$ws.Range("E28").Value = "  -1.08%  "

$ws.Range("D29").Value = "0.0₃0870"
$ws.Range("E29").Value = "  +12.36%  "

$ws.Range("E30").Value = "  +4.69%  "

$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'19.56"
$ws.Range("E32").Value = "  +5.77%  "

$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'156.82"
$ws.Range("E33").Value = "  +4.73%  "

$ws.Range("E34").Value = "  +3.81%  "

$ws.Range("E35").Value = "  +9.61%  "

$ws.Range("E36").Value = "  +8.96%  "

$ws.Range("E37").Value = "  +5.61%  "

$ws.Range("D38").Value = "'311.69"
$ws.Range("E38").Value = "  +16.03%  "

$ws.Range("E39").Value = "  +9.73%  "

$ws.Range("D40").Value = "'0.856"
$ws.Range("E40").Value = "  +2.82%  "

$ws.Range("D41").Value = "'0.846"
$ws.Range("E41").Value = "  +33.65%  "

$ws.Range("E42").Value = "  +7.48%  "

$ws.Range("D43").Value = "'35.20"
$ws.Range("E43").Value = "  +3.39%  "

$ws.Range("D44").Value = "'0.644"
$ws.Range("E44").Value = "  +8.67%  "

$ws.Range("E45").Value = "  +10.36%  "

$ws.Range("D46").Value = "'0.101"
$ws.Range("E46").Value = "  -0.50%  "

$ws.Range("D47").Value = "'20.18"
$ws.Range("E47").Value = "  +16.30%  "

$ws.Range("D48").Value = "'0.996"
$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("D49").Value = "'4.88"
$ws.Range("E49").Value = "  +7.73%  "

$ws.Range("D50").Value = "2.063.33"
$ws.Range("E50").Value = "  +10.63%  "

$ws.Range("E51").Value = "  +3.53%  "
